$d = $word.ActiveDocument

# 1. Update delivery date "05/08" -> "13/01"
$d.Content.Find.Execute("05/08", $true, $false, $false, $false, $false,
                         $true, 1, $false, "13/01", 2)

# 2. Remove "3 ou " from "grupos de 3 ou 4 pessoas"
$d.Content.Find.Execute("3 ou ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
